$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell, forcing text interpretation
# (leading apostrophe) and then stripping the resulting quote-prefix style so
# the cell keeps its original (default) formatting - only the content changes.
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $val
    $r.ClearFormats()
}

# (cell, new value) pairs taken from the authoritative diff.
$updates = @(
    @('D2', '60.855.46'),
    @('E2', '  +0.68%  '),
    @('D3', '2.347.26'),
    @('E3', '  -1.59%  '),
    @('D4', '0.999'),
    @('E4', '  -0.14%  '),
    @('D5', '544.10'),
    @('E5', '  +1.25%  '),
    @('D6', '136.17'),
    @('E6', '  -2.28%  '),
    @('E7', '  +0.00%  '),
    @('E8', '  -7.94%  '),
    @('D9', '2.346.17'),
    @('E9', '  -1.82%  '),
    @('E10', '  +1.45%  '),
    @('E11', '  +1.48%  '),
    @('D12', '5.30'),
    @('E12', '  +0.19%  '),
    @('E13', '  +0.23%  '),
    @('D14', '24.53'),
    @('E14', '  -2.36%  '),
    @('D15', '2.770.07'),
    @('E15', '  -1.64%  '),
    @('D16', '60.661.87'),
    @('E16', '  +0.43%  '),
    @('E17', '  -1.89%  '),
    @('D18', '2.344.48'),
    @('E18', '  -1.72%  '),
    @('D19', '10.59'),
    @('E19', '  +0.81%  '),
    @('D20', '318.77'),
    @('E20', '  +2.19%  '),
    @('E21', '  +1.78%  '),
    @('E22', '  -2.07%  '),
    @('E23', '  +0.02%  '),
    @('E24', '  -2.63%  '),
    @('D25', '63.04'),
    @('E25', '  +0.87%  '),
    @('D26', '8.21'),
    @('E26', '  +8.90%  '),
    @('D27', '7.92'),
    @('E27', '  -0.18%  '),
    @('D28', '495.44'),
    @('E28', '  -0.49%  '),
    @('D29', '1.37'),
    @('E29', '  -2.44%  '),
    @('E30', '  +2.22%  '),
    @('D31', '0.0₃0856'),
    @('E31', '  -4.82%  '),
    @('E32', '  -1.93%  '),
    @('E33', '  -3.22%  '),
    @('E34', '  -0.15%  '),
    @('D35', '4.57'),
    @('E35', '  -0.07%  '),
    @('E36', '  +1.50%  '),
    @('D37', '18.56'),
    @('E37', '  +3.93%  '),
    @('D38', '5.24'),
    @('E38', '  -2.85%  '),
    @('E39', '  +7.80%  '),
    @('D40', '141.47'),
    @('E40', '  +3.94%  '),
    @('E41', '  -0.18%  '),
    @('E42', '  +0.70%  '),
    @('D43', '141.83'),
    @('E43', '  +1.99%  '),
    @('E44', '  +1.66%  '),
    @('E45', '  -4.29%  '),
    @('D46', '0.0514'),
    @('E46', '  +0.60%  '),
    @('B47', 'InjectiveProtocol'),
    @('C47', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'),
    @('D47', '18.88'),
    @('E47', '  -5.23%  '),
    @('B48', 'Stellar'),
    @('C48', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
    @('D48', '0.0899'),
    @('E48', '  -2.00%  '),
    @('B49', 'VeChain'),
    @('C49', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @('D49', '0.0219'),
    @('E49', '  -0.91%  '),
    @('B50', 'Mantle'),
    @('C50', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'),
    @('D50', '0.542'),
    @('E50', '  -5.62%  '),
    @('D51', '16.22'),
    @('E51', '  -1.90%  ')
)

foreach ($u in $updates) {
    Set-TextValue $u[0] $u[1]
}
